$d = $word.ActiveDocument

$replacements = @(
    @{old="25×77=1925"; new="20×70=1400"},
    @{old="91×68=6188"; new="75×60=4500"},
    @{old="42×76=3192"; new="54×99=5346"},
    @{old="44×36=1584"; new="19×63=1197"},
    @{old="16×26=416";  new="65×32=2080"},
    @{old="87×61=5307"; new="46×90=4140"},
    @{old="48×82=3936"; new="22×33=726"},
    @{old="40×39=1560"; new="56×56=3136"},
    @{old="91×99=9009"; new="27×24=648"},
    @{old="45×87=3915"; new="61×48=2928"},
    @{old="95×47=4465"; new="44×63=2772"},
    @{old="20×50=1000"; new="65×73=4745"},
    @{old="87×75=6525"; new="55×25=1375"},
    @{old="47×27=1269"; new="84×19=1596"},
    @{old="94×90=8460"; new="73×61=4453"},
    @{old="57×24=1368"; new="51×63=3213"},
    @{old="26×30=780";  new="71×12=852"},
    @{old="84×28=2352"; new="38×48=1824"},
    @{old="22×89=1958"; new="49×95=4655"},
    @{old="88×27=2376"; new="20×39=780"},
    @{old="40×66=2640"; new="32×63=2016"},
    @{old="64×36=2304"; new="57×92=5244"},
    @{old="50×42=2100"; new="47×26=1222"},
    @{old="27×19=513";  new="31×72=2232"},
    @{old="87×58=5046"; new="29×41=1189"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
